$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.952.44'
$ws.Range('E2').Value = '  +1.90%  '
$ws.Range('D3').Value = '2.525.78'
$ws.Range('E3').Value = '  +1.28%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.48'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '176.71'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.519'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.86%  '
$ws.Range('D9').Value = '2.522.36'
$ws.Range('E9').Value = '  +1.16%  '
$ws.Range('E10').Value = '  +6.05%  '
$ws.Range('E11').Value = '  -1.22%  '
$ws.Range('B12').Value = 'Cardano'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.339'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.54%  '
$ws.Range('B13').Value = 'Toncoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.98'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.20'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.90%  '
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('D16').Value = '68.795.56'
$ws.Range('E16').Value = '  +1.93%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000173'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.68%  '
$ws.Range('D18').Value = '2.512.68'
$ws.Range('E18').Value = '  +3.52%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.12'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.46%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '360.13'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.30%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.53'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.10'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.80%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.80'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.21'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.70'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.79%  '
$ws.Range('E27').Value = '  -3.00%  '
$ws.Range('E28').Value = '  +1.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.997'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '518.98'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.03%  '
$ws.Range('D31').Value = '0.0₃0893'
$ws.Range('E31').Value = '  -1.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.77'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.89%  '
$ws.Range('E33').Value = '  -0.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.78'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.42%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '162.54'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.93%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.119'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.25%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.48'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.70'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.77'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.23%  '
$ws.Range('E41').Value = '  -1.62%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.327'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.83'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.79%  '
$ws.Range('E45').Value = '  -1.32%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '151.55'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.59'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.519'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.07%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0741'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.57%  '
$ws.Range('B50').Value = 'Optimism'
$ws.Range('C50').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.60'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.76%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.580'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.24%  '
